$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update price (D) and volume-change (E) figures, plus the row 44/45 swap
# (PaxDollar <-> TrustWalletToken) to match the latest scrape.

$ws.Range("D2").Value = "30.318.44"
$ws.Range("D3").Value = "1.931.39"
$ws.Range("E4").Value = "  +0.17%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "251.79"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.97%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.7128"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.52%  "
$ws.Range("E7").Value = "  +0.19%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3263"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.62%  "
$ws.Range("E9").Value = "  +3.16%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07204"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +5.53%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.7991"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.90%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08093"
$ws.Range("D12").Style = "Normal"
$ws.Range("D13").Value = "1.929.61"
$ws.Range("E13").Value = "  +0.02%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.428"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.01%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "94.67"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.04%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "14.82"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.75%  "
$ws.Range("D17").Value = "30.315.29"
$ws.Range("E17").Value = "  +0.17%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "251.98"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -3.68%  "
$ws.Range("E19").Value = "  +2.19%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "5.787"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.64%  "
$ws.Range("D21").Value = "2.180.27"
$ws.Range("E21").Value = "  -0.09%  "
$ws.Range("E22").Value = "  +0.16%  "
$ws.Range("E23").Value = "  +0.22%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.921"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.18%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.712"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.03%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "164.79"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +3.03%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "19.22"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.00%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.318"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.62%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.1283"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -5.22%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.359"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.50%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.543"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.44%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.433"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.73%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.210"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.09%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.05207"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +2.32%  "
$ws.Range("E35").Value = "  +5.09%  "
$ws.Range("E36").Value = "  +0.63%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.760"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.35%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01964"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.19%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.800"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.38%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "78.81"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.89%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.422"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.43%  "
$ws.Range("E42").Value = "  +1.05%  "
$ws.Range("E43").Value = "  +0.80%  "
$ws.Range("B44").Value = "TrustWalletToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.8415"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.76%  "
$ws.Range("B45").Value = "PaxDollar"
$ws.Range("C45").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.001"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.14%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "101.79"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.65%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.851"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.02%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.423"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.54%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "36.72"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.01%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06089"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.84%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.4178"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.28%  "
